$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantIntTable")

# Fix Gacha3Events value: 10 -> 9 (row 11, column B)
$ws.Range("B11").Value = 9

# Insert a new row at row 14 (shifts existing rows 14-19 down to 15-20)
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new constant
$ws.Range("A14").Value = "Gacha3BrokenEnergys"
$ws.Range("B14").Value = 3
